$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# sheet1
$ws1.Range("F3").Value = 3845
$ws1.Range("F6").Value = 3841
$ws1.Range("F8").Value = 207
$ws1.Range("F9").Value = 59
$ws1.Range("F10").Value = 8740
$ws1.Range("F12").Value = 136
$ws1.Range("F13").Value = 306
$ws1.Range("F14").Value = 346
$ws1.Range("F15").Value = 135
$ws1.Range("F16").Value = 105
$ws1.Range("F18").Value = 371
$ws1.Range("F19").Value = 11083
$ws1.Range("F28").Value = 194
$ws1.Range("F38").Value = 908
$ws1.Range("F40").Value = 287
$ws1.Range("F42").Value = 1253
$ws1.Range("F44").Value = 759
$ws1.Range("F46").Value = 351
$ws1.Range("F47").Value = 50

# sheet2
$ws2.Range("F9").Value = 10
$ws2.Range("F11").Value = 33
$ws2.Range("F15").Value = 37
$ws2.Range("F16").Value = 13
$ws2.Range("F22").Value = 63

# sheet3
$ws3.Range("F3").Value = 36

# sheet4
$ws4.Range("F4").Value = 3845
$ws4.Range("F6").Value = 3841
$ws4.Range("F10").Value = 207
$ws4.Range("F11").Value = 59
$ws4.Range("F12").Value = 8740
$ws4.Range("F13").Value = 33
$ws4.Range("F15").Value = 136
$ws4.Range("F16").Value = 306
$ws4.Range("F17").Value = 346
$ws4.Range("F18").Value = 135
$ws4.Range("F19").Value = 105
$ws4.Range("F20").Value = 371
$ws4.Range("F21").Value = 11083
$ws4.Range("F26").Value = 194
$ws4.Range("F27").Value = 37
$ws4.Range("F34").Value = 909
$ws4.Range("F38").Value = 287
$ws4.Range("F40").Value = 63
$ws4.Range("F41").Value = 1253
$ws4.Range("F43").Value = 759
$ws4.Range("F45").Value = 351
$ws4.Range("F47").Value = 50
